$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "51.583.12"
Set-TextValue $ws.Range("E2") "  +1.11%  "

Set-TextValue $ws.Range("D3") "2.988.71"
Set-TextValue $ws.Range("E3") "  +2.77%  "

Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  +0.00%  "

Set-TextValue $ws.Range("D5") "381.90"
Set-TextValue $ws.Range("E5") "  +4.63%  "

Set-TextValue $ws.Range("D6") "106.04"
Set-TextValue $ws.Range("E6") "  +2.83%  "

Set-TextValue $ws.Range("E7") "  +1.10%  "

Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  +0.07%  "

Set-TextValue $ws.Range("D9") "0.599"
Set-TextValue $ws.Range("E9") "  +1.78%  "

Set-TextValue $ws.Range("D10") "37.47"
Set-TextValue $ws.Range("E10") "  +1.53%  "

Set-TextValue $ws.Range("E11") "  +0.66%  "

Set-TextValue $ws.Range("D12") "0.0845"
Set-TextValue $ws.Range("E12") "  +1.28%  "

Set-TextValue $ws.Range("D13") "18.69"
Set-TextValue $ws.Range("E13") "  +1.67%  "

Set-TextValue $ws.Range("D14") "3.458.92"
Set-TextValue $ws.Range("E14") "  +2.75%  "

Set-TextValue $ws.Range("E15") "  +2.45%  "

Set-TextValue $ws.Range("D16") "2.991.34"
Set-TextValue $ws.Range("E16") "  +3.10%  "

Set-TextValue $ws.Range("D17") "0.972"
Set-TextValue $ws.Range("E17") "  +2.21%  "

Set-TextValue $ws.Range("D18") "51.649.56"
Set-TextValue $ws.Range("E18") "  +1.27%  "

Set-TextValue $ws.Range("E19") "  +3.53%  "

Set-TextValue $ws.Range("D20") "7.42"
Set-TextValue $ws.Range("E20") "  +2.46%  "

Set-TextValue $ws.Range("E21") "  +0.54%  "

Set-TextValue $ws.Range("E22") "  +2.08%  "

Set-TextValue $ws.Range("D23") "69.25"
Set-TextValue $ws.Range("E23") "  +1.71%  "

Set-TextValue $ws.Range("D24") "263.64"
Set-TextValue $ws.Range("E24") "  +1.35%  "

Set-TextValue $ws.Range("D25") "2.81"
Set-TextValue $ws.Range("E25") "  +4.54%  "

Set-TextValue $ws.Range("E26") "  -1.19%  "

Set-TextValue $ws.Range("D27") "7.26"
Set-TextValue $ws.Range("E27") "  +19.01%  "

Set-TextValue $ws.Range("D28") "7.54"
Set-TextValue $ws.Range("E28") "  +3.24%  "

Set-TextValue $ws.Range("D30") "26.09"
Set-TextValue $ws.Range("E30") "  +0.74%  "

Set-TextValue $ws.Range("D31") "0.109"
Set-TextValue $ws.Range("E31") "  +4.87%  "

Set-TextValue $ws.Range("D32") "9.91"
Set-TextValue $ws.Range("E32") "  -0.05%  "

Set-TextValue $ws.Range("D33") "35.12"
Set-TextValue $ws.Range("E33") "  +0.45%  "

Set-TextValue $ws.Range("E34") "  +10.03%  "

Set-TextValue $ws.Range("D35") "2.09"
Set-TextValue $ws.Range("E35") "  -2.11%  "

Set-TextValue $ws.Range("D36") "51.36"
Set-TextValue $ws.Range("E36") "  +1.58%  "

Set-TextValue $ws.Range("E37") "  +0.13%  "

Set-TextValue $ws.Range("E38") "  -0.70%  "

Set-TextValue $ws.Range("D39") "17.52"
Set-TextValue $ws.Range("E39") "  +3.77%  "

Set-TextValue $ws.Range("D40") "2.62"
Set-TextValue $ws.Range("E40") "  -5.88%  "

Set-TextValue $ws.Range("E41") "  -0.15%  "

Set-TextValue $ws.Range("E42") "  +2.83%  "

Set-TextValue $ws.Range("D43") "123.81"
Set-TextValue $ws.Range("E43") "  +5.29%  "

Set-TextValue $ws.Range("D44") "22.43"
Set-TextValue $ws.Range("E44") "  -0.24%  "

Set-TextValue $ws.Range("D45") "0.283"
Set-TextValue $ws.Range("E45") "  +20.91%  "

Set-TextValue $ws.Range("E46") "  -0.56%  "

Set-TextValue $ws.Range("D47") "2.39"
Set-TextValue $ws.Range("E47") "  +5.56%  "

Set-TextValue $ws.Range("D48") "2.055.81"
Set-TextValue $ws.Range("E48") "  -0.34%  "

Set-TextValue $ws.Range("E49") "  +2.11%  "

Set-TextValue $ws.Range("D50") "0.0356"
Set-TextValue $ws.Range("E50") "  +11.18%  "

Set-TextValue $ws.Range("E51") "  +3.79%  "
